$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.076.71'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '1.835.83'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6248'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.90%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07558'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2918'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = '1.839.06'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.951'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6656'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009914'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +12.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.038'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '29.113.58'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.186'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.455'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1368'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.493'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.076'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.035'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.202'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05200'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.856'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7409'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.704'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.87%  '
$ws.Range("D37").Value = '1.250.48'
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.766'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01783'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.377'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8934'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.83%  '
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("D44").Value = '1.981.82'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5120'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4016'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.839'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.652'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05755'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.91%  '
